# Edit script: applies the 2023-10-31 liga-3 scrape update.
# 1) Row 27 <-> Row 31 had their match data (cols F:V) swapped back to the
#    correct fixtures (Canelas 2010 vs Felgueiras / Academica vs Amora).
# 2) Row 77 <-> Row 79 had their match data (cols F:V) swapped back to the
#    correct fixtures (Varzim vs Sanjoanense / Amora vs 1o Dezembro).
# 3) Ten new matches (rows 81-90) were appended for fixtures played
#    2023-10-26 through 2023-10-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, [int]$row1, [int]$row2)
    # columns F..V (6..22) hold the match-specific data; A..E identify the
    # fixture slot itself and stay untouched.
    # NB: use Value2 (not Value) for reads -- Value's getter is unreliable
    # in this host and can yield the property descriptor instead of data.
    for ($c = 6; $c -le 22; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)
        $tmp = $cell1.Value2
        $cell1.Value2 = $cell2.Value2
        $cell2.Value2 = $tmp
    }
}

Swap-RowData $ws 27 31
Swap-RowData $ws 77 79

# Append rows 81-90, matching the look & feel (styles) of the existing data
# rows by copying the format from the last populated row (80) first.
$ws.Range("A80:V80").Copy()
$ws.Range("A81:V90").PasteSpecial(-4122)

$newRows = @(
    ,@(80, "portugal", "liga-3", "2023-2024", 45226.85416666666, "Alverca", 1, "Amora", 0, 1.65, "26/10/2023 13:54", 1.49, "27/10/2023 20:29", 3.67, "26/10/2023 13:54", 4.26, "27/10/2023 20:29", 4.98, "26/10/2023 13:54", 6.88, "27/10/2023 20:29", "https://www.betexplorer.com/football/portugal/liga-3/alverca-amora/COcxVctJ/")
    ,@(81, "portugal", "liga-3", "2023-2024", 45227.70833333334, "Covilha", 4, "Caldas", 3, 1.99, "27/10/2023 13:43", 2.35, "28/10/2023 16:51", 3.27, "27/10/2023 13:43", 3.36, "28/10/2023 16:51", 3.76, "27/10/2023 13:43", 3.11, "28/10/2023 16:51", "https://www.betexplorer.com/football/portugal/liga-3/covilha-caldas-sc/rq6LYusf/")
    ,@(82, "portugal", "liga-3", "2023-2024", 45228.66666666666, "1º Dezembro", 0, "Oliveira Hospital", 2, 3.47, "28/10/2023 18:13", 2.92, "29/10/2023 12:06", 3.24, "28/10/2023 18:13", 3.36, "29/10/2023 12:06", 2.18, "28/10/2023 18:13", 2.48, "29/10/2023 12:06", "https://www.betexplorer.com/football/portugal/liga-3/1-dezembro-oliveira-hospital/YXbYVwRC/")
    ,@(83, "portugal", "liga-3", "2023-2024", 45228.66666666666, "Atletico CP", 2, "Academica", 1, 2.04, "28/10/2023 18:13", 2.4, "29/10/2023 12:34", 3.31, "28/10/2023 18:13", 3.2, "29/10/2023 12:34", 3.57, "28/10/2023 18:13", 3.18, "29/10/2023 12:34", "https://www.betexplorer.com/football/portugal/liga-3/atletico-cp-academica/x65PXad0/")
    ,@(84, "portugal", "liga-3", "2023-2024", 45228.66666666666, "Pero Pinheiro", 3, "Sporting CP B", 3, 3.92, "28/10/2023 18:13", 6.79, "29/10/2023 15:42", 3.33, "28/10/2023 18:13", 4.37, "29/10/2023 15:42", 2, "28/10/2023 18:13", 1.48, "29/10/2023 15:41", "https://www.betexplorer.com/football/portugal/liga-3/pero-pinheiro-sporting-lisbon/j99TWJB6/")
    ,@(85, "portugal", "liga-3", "2023-2024", 45228.66666666666, "Braga B", 0, "SC Vianense", 3, 1.5, "28/10/2023 18:13", 1.56, "29/10/2023 10:50", 4.24, "28/10/2023 18:13", 4.3, "29/10/2023 10:50", 6.35, "28/10/2023 18:13", 5.61, "29/10/2023 10:50", "https://www.betexplorer.com/football/portugal/liga-3/braga-sc-vianense/b3WFOt1d/")
    ,@(86, "portugal", "liga-3", "2023-2024", 45228.66666666666, "Canelas 2010", 3, "AD Fafe", 1, 2.14, "28/10/2023 18:13", 2.58, "29/10/2023 12:25", 3.27, "28/10/2023 18:13", 3.31, "29/10/2023 08:22", 3.55, "28/10/2023 18:13", 2.83, "29/10/2023 12:25", "https://www.betexplorer.com/football/portugal/liga-3/canelas-2010-ad-fafe/QkY7Q2Wq/")
    ,@(87, "portugal", "liga-3", "2023-2024", 45228.66666666666, "Sanjoanense", 1, "Lusitania FC", 1, 4.15, "28/10/2023 18:13", 6.18, "29/10/2023 13:08", 3.5, "28/10/2023 18:13", 4.15, "29/10/2023 15:49", 1.88, "28/10/2023 18:13", 1.55, "29/10/2023 13:08", "https://www.betexplorer.com/football/portugal/liga-3/sanjoanense-lusitania-fc/ED5vrV1H/")
    ,@(88, "portugal", "liga-3", "2023-2024", 45228.75, "Anadia", 3, "Trofense", 3, 2.12, "28/10/2023 18:13", 2.26, "29/10/2023 17:51", 3.27, "28/10/2023 18:13", 3.32, "29/10/2023 17:51", 3.61, "28/10/2023 18:13", 3.32, "29/10/2023 17:51", "https://www.betexplorer.com/football/portugal/liga-3/anadia-trofense/WlGZrBnB/")
    ,@(89, "portugal", "liga-3", "2023-2024", 45228.77083333334, "Felgueiras", 2, "Varzim", 0, 1.84, "28/10/2023 18:13", 1.71, "29/10/2023 18:21", 3.47, "28/10/2023 18:13", 3.7, "29/10/2023 18:21", 4.46, "28/10/2023 18:13", 5.13, "29/10/2023 18:21", "https://www.betexplorer.com/football/portugal/liga-3/fc-felgueiras-varzim/KCXBPMoj/")
)

$r = 81
foreach ($row in $newRows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
